$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New listing: Count Basiestraat 52 - add as row 18, following the same
# layout/formatting as the row above it (row 17).
$ws.Range("A17:F17").Copy($ws.Range("A18:F18"))

$ws.Range("A18").Value = "Count Basiestraat 52"
$ws.Range("B18").Value = 229500
$ws.Range("C18").Value = "1311 PD Almere"
$ws.Range("D18").Value = "Email Sent"
$ws.Range("E18").Value = "No Updates yet"
$ws.Range("F18").Value = "https://www.funda.nl/koop/almere/huis-40693517-count-basiestraat-52/"

# Match the author's final selection in the saved workbook.
$ws.Range("C21").Select()
